# Update iServ stats for 2025-10 (row 23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = 6321
$ws.Range("C23").Value = 997
$ws.Range("D23").Value = 5885516
$ws.Range("E23").Value = 931.1052048726467
$ws.Range("F23").Value = 8.45916266300617
$ws.Range("G23").Value = 3.746097814776284
$ws.Range("H23").Value = 26.11154311211716
